# This script updates various item description / name strings in the
# "Worksheet" sheet (column C holds the English text) to match the
# translation corrections described in the commit:
#   "Updated Items, Troops, and Skills.xlsx with various translations."
#
# Most of the changes replace the mistranslated term "strength" (referring
# to HP/physical-strength stat) with "HP", and a couple of "mental
# strength" occurrences are replaced with "Willpower". Also "Sinapis" /
# "Sina" are corrected to "Shinapis" / "Shina".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$replacements = @{
    "C4"   = "Tastes lightly of orange.`nRestores 30% + 20 of your HP."
    "C32"  = "Tastes mildly like oranges.`nRestores 30% +5 of HP."
    "C35"  = "Tastes like oranges.`nRestores 50% +5 of HP."
    "C38"  = "Tastes strongly of oranges.`nRestores 70% of HP."
    "C41"  = "Overwhelmingly orange flavored.`nFully restores HP."
    "C44"  = "Ground and flavored medicinal herbs which are easy to drink.`nRestores 30% vitality and HP."
    "C47"  = "An adventurer's best friend which energizes both body and soul.`nRestores 50% vitality and HP."
    "C50"  = "Ambient magic which has been condensed into liquid form.`n70% recovery of vitality and HP."
    "C53"  = "An alchemical product created through unknown means.`nFull recovery of vitality and HP, and a large increase to TP"
    "C64"  = "A potion that temporarily increases defense and Willpower.`nDefensive power up, Mental power up, small increase to TP."
    "C69"  = "A card with the power to twist fate.`nFull restoration of vigor and HP, and all abilities UP."
    "C74"  = "Water of Shinapis"
    "C75"  = "A panacea for all that ails.`nRestores all vitality and HP, extra large increase to TP."
    "C85"  = "Nourishes the body and restores HP.`nRestores 20 HP, and removes weakness."
    "C87"  = "Sedates pain and restores one to full HP.`nFull restoration of HP and removes weakness."
    "C93"  = "Sweet pastry made with date bean paste.`nRecovers 5 HP"
    "C101" = "Peach juice, said to be favored by the gods.`nDefenses UP, Willpower UP, maximum vigor UP."
    "C110" = "A berry that makes you mentally strong when you eat it.`nWillpower +1 (permanent)."
    "C118" = "Berry that make it easier to gain muscle mass when eaten.`nMaximum HP +5 (permanent)."
    "C122" = "HP herb"
    "C287" = "Memories of Defeat: Shina"
    "C288" = "Shina's memories of defeat.`nHow many you've got shows how much you've been squeezed by them"
}

foreach ($addr in $replacements.Keys) {
    $ws.Range($addr).Value = $replacements[$addr]
}
